$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "331.49"
Set-TextValue $ws.Range("E2") "0.27%"
Set-TextValue $ws.Range("D3") "41.72"
Set-TextValue $ws.Range("E3") "6.62%"
Set-TextValue $ws.Range("D4") "5.715"
Set-TextValue $ws.Range("E4") "0.84%"
Set-TextValue $ws.Range("D5") "0.08345"
Set-TextValue $ws.Range("E5") "4.00%"
Set-TextValue $ws.Range("D6") "2.031"
Set-TextValue $ws.Range("E6") "4.58%"
Set-TextValue $ws.Range("D7") "8.796"
Set-TextValue $ws.Range("E7") "2.27%"
Set-TextValue $ws.Range("E8") "1.37%"
Set-TextValue $ws.Range("E9") "2.28%"
Set-TextValue $ws.Range("D10") "0.9258"
Set-TextValue $ws.Range("E10") "0.73%"
Set-TextValue $ws.Range("D11") "0.1289"
Set-TextValue $ws.Range("E11") "4.15%"
Set-TextValue $ws.Range("D12") "0.1962"
Set-TextValue $ws.Range("E12") "1.28%"
Set-TextValue $ws.Range("D13") "0.09412"
Set-TextValue $ws.Range("E13") "2.37%"
Set-TextValue $ws.Range("D14") "0.03908"
Set-TextValue $ws.Range("E14") "11.40%"
Set-TextValue $ws.Range("D15") "0.1061"
Set-TextValue $ws.Range("E15") "1.09%"
Set-TextValue $ws.Range("D16") "0.001308"
Set-TextValue $ws.Range("E16") "-0.51%"
Set-TextValue $ws.Range("D17") "0.006138"
Set-TextValue $ws.Range("E17") "-3.01%"
Set-TextValue $ws.Range("D18") "3.440"
Set-TextValue $ws.Range("E18") "2.31%"
Set-TextValue $ws.Range("E19") "2.28%"
Set-TextValue $ws.Range("D20") "8.228"
Set-TextValue $ws.Range("E20") "-5.40%"
Set-TextValue $ws.Range("D21") "0.1373"
Set-TextValue $ws.Range("D22") "0.2413"
Set-TextValue $ws.Range("E22") "-10.54%"
Set-TextValue $ws.Range("D23") "0.04416"
Set-TextValue $ws.Range("E23") "-0.54%"
Set-TextValue $ws.Range("D24") "0.001249"
Set-TextValue $ws.Range("E24") "-0.47%"
Set-TextValue $ws.Range("D25") "0.004380"
Set-TextValue $ws.Range("E25") "-2.61%"
Set-TextValue $ws.Range("D26") "0.0001200"
Set-TextValue $ws.Range("E26") "-0.17%"
Set-TextValue $ws.Range("D39") "0.02795"
Set-TextValue $ws.Range("E39") "9.55%"
Set-TextValue $ws.Range("D40") "0.05566"
Set-TextValue $ws.Range("E40") "2.38%"
Set-TextValue $ws.Range("D41") "0.007795"
Set-TextValue $ws.Range("E41") "3.40%"
Set-TextValue $ws.Range("E42") "2.38%"
Set-TextValue $ws.Range("D43") "0.008947"
Set-TextValue $ws.Range("E43") "-9.79%"
Set-TextValue $ws.Range("D44") "0.002140"
Set-TextValue $ws.Range("E44") "1.38%"
Set-TextValue $ws.Range("D45") "0.01053"
Set-TextValue $ws.Range("E45") "-8.73%"
Set-TextValue $ws.Range("D46") "0.00007046"
Set-TextValue $ws.Range("E46") "3.57%"
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.17%"
Set-TextValue $ws.Range("D48") "0.003524"
Set-TextValue $ws.Range("E48") "15.30%"
Set-TextValue $ws.Range("D49") "0.002280"
Set-TextValue $ws.Range("E49") "-0.12%"
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "-0.17%"
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.17%"
